# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# on the Leve profit sheets with the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3389.2354
$ws.Range("I62").Value = 3366.6667
$ws.Range("J62").Value = 3443.4
$ws.Range("K62").Value = 3366.6667
$ws.Range("L62").Value = 3443.4
$ws.Range("M62").Value = -2742.6667
$ws.Range("N62").Value = -4691.4

$ws.Range("H64").Value = 2923.35
$ws.Range("I64").Value = 2807.25
$ws.Range("J64").Value = 3000.75
$ws.Range("K64").Value = 2807.25
$ws.Range("L64").Value = 3000.75
$ws.Range("M64").Value = -2559.25
$ws.Range("N64").Value = -3496.75

$ws.Range("H65").Value = 3389.2354
$ws.Range("I65").Value = 3366.6667
$ws.Range("J65").Value = 3443.4
$ws.Range("K65").Value = 16833.3335
$ws.Range("L65").Value = 17217
$ws.Range("M65").Value = -13713.3335
$ws.Range("N65").Value = -23457

$ws.Range("H67").Value = 2923.35
$ws.Range("I67").Value = 2807.25
$ws.Range("J67").Value = 3000.75
$ws.Range("K67").Value = 2807.25
$ws.Range("L67").Value = 3000.75
$ws.Range("M67").Value = -1949.25
$ws.Range("N67").Value = -4716.75

$ws.Range("H76").Value = 3327.3333
$ws.Range("I76").Value = 2921.7273
$ws.Range("J76").Value = 3964.7144
$ws.Range("K76").Value = 2921.7273
$ws.Range("L76").Value = 3964.7144
$ws.Range("M76").Value = -2606.7273
$ws.Range("N76").Value = -4594.7144

$ws.Range("H79").Value = 3327.3333
$ws.Range("I79").Value = 2921.7273
$ws.Range("J79").Value = 3964.7144
$ws.Range("K79").Value = 2921.7273
$ws.Range("L79").Value = 3964.7144
$ws.Range("M79").Value = -1829.7273
$ws.Range("N79").Value = -6148.7144

$ws.Range("H92").Value = 864.5714
$ws.Range("I92").Value = 399
$ws.Range("J92").Value = 1213.75
$ws.Range("K92").Value = 399
$ws.Range("L92").Value = 1213.75
$ws.Range("M92").Value = 849
$ws.Range("N92").Value = -3709.75

$ws.Range("H99").Value = 3830.3635
$ws.Range("J99").Value = 40000
$ws.Range("L99").Value = 120000
$ws.Range("N99").Value = -122996

$ws.Range("H106").Value = 1153.3334
$ws.Range("I106").Value = 1153.3334
$ws.Range("K106").Value = 1153.3334
$ws.Range("M106").Value = -522.3334

$ws.Range("H129").Value = 6880.972
$ws.Range("I129").Value = 438
$ws.Range("K129").Value = 1314
$ws.Range("M129").Value = 3686

$ws.Range("H138").Value = 2432.6365
$ws.Range("I138").Value = 1258.2258
$ws.Range("K138").Value = 3774.6774
$ws.Range("M138").Value = 1365.3226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4485
$ws.Range("I94").Value = 475.48648
$ws.Range("J94").Value = 15081.571
$ws.Range("K94").Value = 475.48648
$ws.Range("L94").Value = 15081.571
$ws.Range("M94").Value = -24.48647999999997
$ws.Range("N94").Value = -15983.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 14357.143
$ws.Range("J80").Value = 14357.143
$ws.Range("L80").Value = 14357.143
$ws.Range("N80").Value = -16603.143

$ws.Range("H83").Value = 14357.143
$ws.Range("J83").Value = 14357.143
$ws.Range("L83").Value = 43071.429
$ws.Range("N83").Value = -54303.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 831.69696
$ws.Range("I113").Value = 434.07144
$ws.Range("J113").Value = 1124.6842
$ws.Range("K113").Value = 1302.21432
$ws.Range("L113").Value = 3374.0526
$ws.Range("M113").Value = 867.78568
$ws.Range("N113").Value = -7714.0526

$ws.Range("H117").Value = 4175
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 4175
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 12525
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -19409

$ws.Range("H129").Value = 1008.58826
$ws.Range("I129").Value = 1000
$ws.Range("J129").Value = 1009.125
$ws.Range("K129").Value = 3000
$ws.Range("L129").Value = 3027.375
$ws.Range("M129").Value = 2000
$ws.Range("N129").Value = -13027.375

$ws.Range("H131").Value = 1011.2
$ws.Range("I131").Value = 445.5
$ws.Range("K131").Value = 1336.5
$ws.Range("M131").Value = 3703.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 25800.2
$ws.Range("I34").Value = 12500
$ws.Range("J34").Value = 34667
$ws.Range("K34").Value = 12500
$ws.Range("L34").Value = 34667
$ws.Range("M34").Value = -12232
$ws.Range("N34").Value = -35203

$ws.Range("H76").Value = 25800.2
$ws.Range("I76").Value = 12500
$ws.Range("J76").Value = 34667
$ws.Range("K76").Value = 12500
$ws.Range("L76").Value = 34667
$ws.Range("M76").Value = -12185
$ws.Range("N76").Value = -35297

$ws.Range("H79").Value = 25800.2
$ws.Range("I79").Value = 12500
$ws.Range("J79").Value = 34667
$ws.Range("K79").Value = 12500
$ws.Range("L79").Value = 34667
$ws.Range("M79").Value = -11408
$ws.Range("N79").Value = -36851

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1825
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1825
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1825
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2415

$ws.Range("H27").Value = 1825
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1825
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1825
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -2039

$ws.Range("H87").Value = 29500.5
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 29500.5
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 29500.5
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -31746.5

$ws.Range("H88").Value = 30001
$ws.Range("J88").Value = 30001
$ws.Range("L88").Value = 30001
$ws.Range("N88").Value = -30857

$ws.Range("H90").Value = 29500.5
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 29500.5
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 88501.5
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -99733.5

$ws.Range("H91").Value = 30001
$ws.Range("J91").Value = 30001
$ws.Range("L91").Value = 30001
$ws.Range("N91").Value = -32965

$ws.Range("H93").Value = 1614.1111
$ws.Range("I93").Value = 1548.8125
$ws.Range("J93").Value = 1709.091
$ws.Range("K93").Value = 1548.8125
$ws.Range("L93").Value = 1709.091
$ws.Range("M93").Value = -300.8125
$ws.Range("N93").Value = -4205.091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 24222.25
$ws.Range("J75").Value = 24222.25
$ws.Range("L75").Value = 24222.25
$ws.Range("N75").Value = -26094.25

$ws.Range("H78").Value = 24222.25
$ws.Range("J78").Value = 24222.25
$ws.Range("L78").Value = 72666.75
$ws.Range("N78").Value = -82026.75

$ws.Range("H109").Value = 8580
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 8580
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 8580
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -11354

Write-Output "done"
